$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the phone column (E2:E23) to hold the text value "+84947947990"
# instead of the bare numeric digits. Forcing text via a leading apostrophe
# so Excel stores it as text (quote-prefixed) rather than re-parsing it as
# a number.
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 5).Value = "'+84947947990"
}

# Widen column E to fit the new text and drop the old autosized bestFit.
# (ColumnWidth is specified in characters; 18.8333... is the COM-level
# value that serializes to the OOXML column width of ~19.71 used by the
# workbook author.)
$ws.Columns.Item(5).ColumnWidth = 18.8333333333333

# Restore the view: selecting I34 clears the old scrolled-down
# topLeftCell/selection state and leaves just the new single-cell selection.
$ws.Range("I34").Select()
